$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '307.88'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '0.13%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '40.96'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '2.36%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.104'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-0.25%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07620'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-1.45%'
$ws.Range("B6").Value = 'GateToken'
$ws.Range("C6").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '4.251'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '0.37%'
$ws.Range("B7").Value = 'FTXToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.604'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '-0.67%'
$ws.Range("B8").Value = 'BTSEToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '2.471'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '2.12%'
$ws.Range("B9").Value = 'MXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9034'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '1.94%'
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1123'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '12.89%'
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1783'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '2.39%'
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.09156'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '1.08%'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.04225'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-5.19%'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.1053'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.12%'
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001255'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-0.12%'
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.005752'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-2.03%'
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.349'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-0.10%'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.639'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-6.85%'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '1.23%'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '1.49%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.04076'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-0.92%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.001245'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '3.69%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.004109'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '0.97%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0001302'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-0.02%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0003749'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02376'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '1.17%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05180'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '-0.44%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.007770'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-1.90%'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-1.51%'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '12.42%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.001953'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '0.01%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.3079'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-7.50%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00007005'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '6.76%'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-0.03%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.03115'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '765.21%'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.03%'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '-0.03%'
